$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their literal text representation
# (values like "307.69" or "3.41%" would otherwise be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "307.69"
$ws.Range("E2").Value = "3.41%"
$ws.Range("D3").Value = "44.19"
$ws.Range("E3").Value = "6.38%"
$ws.Range("D4").Value = "5.073"
$ws.Range("E4").Value = "0.91%"
$ws.Range("D5").Value = "0.07979"
$ws.Range("E5").Value = "5.76%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "4.428"
$ws.Range("E6").Value = "1.22%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.599"
$ws.Range("E7").Value = "0.42%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "1.071"
$ws.Range("E8").Value = "15.29%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1280"
$ws.Range("E9").Value = "6.84%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1889"
$ws.Range("E10").Value = "2.41%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.09184"
$ws.Range("E11").Value = "3.03%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.04202"
$ws.Range("E12").Value = "3.83%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.1036"
$ws.Range("E13").Value = "-1.73%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001291"
$ws.Range("E14").Value = "0.85%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005710"
$ws.Range("E15").Value = "-0.93%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "0.007409"
$ws.Range("E16").Value = "1,889.69%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.338"
$ws.Range("E17").Value = "-0.05%"
$ws.Range("D18").Value = "2.394"
$ws.Range("E18").Value = "-1.22%"
$ws.Range("D19").Value = "0.3348"
$ws.Range("E19").Value = "1.14%"
$ws.Range("D20").Value = "8.045"
$ws.Range("E20").Value = "0.92%"
$ws.Range("D21").Value = "0.1362"
$ws.Range("E21").Value = "-4.01%"
$ws.Range("D23").Value = "0.04140"
$ws.Range("E23").Value = "2.21%"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").Value = "-0.08%"
$ws.Range("D25").Value = "0.004571"
$ws.Range("E25").Value = "8.95%"
$ws.Range("D26").Value = "0.0001333"
$ws.Range("E26").Value = "8.39%"
$ws.Range("D38").Value = "0.02674"
$ws.Range("E38").Value = "10.68%"
$ws.Range("D39").Value = "0.05364"
$ws.Range("E39").Value = "2.78%"
$ws.Range("D40").Value = "0.005597"
$ws.Range("E40").Value = "-15.20%"
$ws.Range("D41").Value = "0.007725"
$ws.Range("E41").Value = "-1.19%"
$ws.Range("D42").Value = "0.1403"
$ws.Range("E42").Value = "5.50%"
$ws.Range("D43").Value = "0.007224"
$ws.Range("E43").Value = "-4.52%"
$ws.Range("D44").Value = "0.008401"
$ws.Range("E44").Value = "7.30%"
$ws.Range("D45").Value = "0.3067"
$ws.Range("E45").Value = "-4.82%"
$ws.Range("D46").Value = "0.00006558"
$ws.Range("E46").Value = "-1.51%"
$ws.Range("E47").Value = "-1.26%"
$ws.Range("D48").Value = "0.05119"
$ws.Range("E48").Value = "10.65%"
$ws.Range("D49").Value = "0.003923"
$ws.Range("E49").Value = "-6.58%"
$ws.Range("D50").Value = "0.00002073"
$ws.Range("E50").Value = "-1.26%"
$ws.Range("E51").Value = "-1.26%"
